# Update the timestamp portion of the test email addresses from
# "20251109_003734" to "20251109_004215" everywhere they appear in the workbook.

$wb = $excel.ActiveWorkbook

$oldStamp = "20251109_003734"
$newStamp = "20251109_004215"

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        $val = $cell.Value2
        if ($val -ne $null -and $val -is [string] -and $val.Contains($oldStamp)) {
            $cell.Value = $val.Replace($oldStamp, $newStamp)
        }
    }
}

$wb.Save()
